# NSMB - World 5 begin
# Applies the "World 5" restructuring to the V4 sheet (sheet1 / ActiveSheet):
#  - adds K/L helper cells (style 13) to rows 15-20
#  - records new "current run" Place values in column B for rows 30, 32, 33, 36, 37
#    (their D/F diff formulas auto-recalculate)
#  - for rows 36-57, the old "Place"/"Best" pair (columns B/C) is moved into new
#    "Me"/"Was0x" columns (G/H, header added in row 35), and B/C is cleared for
#    rows 38-57 (their diff formulas fall back to 0)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reusable style donor cells (never modified below) so new cells pick up the
# SAME existing style index instead of Excel minting a fresh one.
$style13 = $ws.Range("K8")    # plain helper-column style used by K/L columns
$style16 = $ws.Range("A17")   # thin-right-border style used by B..H data columns
$style22 = $ws.Range("B41")   # section-header fill+border style used by B..H

function Set-Style($range, $donor) {
    $donor.Copy()
    $range.PasteSpecial(-4122)   # xlPasteFormats
}

# ---- Rows 15-20: add K/L "placeholder" cells ----
Set-Style $ws.Range("L15") $style13
Set-Style $ws.Range("L16") $style13
Set-Style $ws.Range("K17:L17") $style13
Set-Style $ws.Range("K18:L18") $style13
Set-Style $ws.Range("K19:L19") $style13
Set-Style $ws.Range("L20") $style13

# ---- Rows 30, 32, 33: new "current run" Place values (column B) ----
# D/F already hold the IF(B>0, ...) formulas, so setting B recalculates them.
$ws.Range("B30").Value = 11140
$ws.Range("B32").Value = 11592
$ws.Range("B33").Value = 12540
Set-Style $ws.Range("B33") $style16

# ---- Row 35: new "Me" / "Was0x" header over columns G/H ----
Set-Style $ws.Range("G35:H35") $style22
$ws.Range("G35").Value = "Me"
$ws.Range("H35").Value = "Was0x"
# header cells carry no explicit style in the diff (same as row 35's other
# blank cells before B..F got theirs) - clear the donor formatting again.
$ws.Range("G35").ClearFormats()
$ws.Range("H35").ClearFormats()
$ws.Range("G35").Value = "Me"
$ws.Range("H35").Value = "Was0x"

# ---- Rows 36-57: move old Place/Best (B/C) into the new Me/Was0x columns (G/H) ----
$oldValues = @{
    36 = @(15202, 15292)
    37 = @(16484, 16574)
    38 = @(17011, 17101)
    39 = @(17525, 17615)
    40 = @(18058, 18298)
    42 = @(18356, 18596)
    43 = @(18742, 18982)
    44 = @(18974, 19214)
    45 = @(19152, 19392)
    46 = @(19206, 19447)
    47 = @(19281, 19522)
    48 = @(19350, 19591)
    49 = @(19374, 19615)
    50 = @(19461, 19703)
    51 = @(19706, 19949)
    52 = @(20114, 20359)
    53 = @(20117, 20363)
    54 = @(20257, 20512)
    55 = @(20537, 20832)
    56 = @(21051, 21346)
    57 = @(21617, 22007)
}

foreach ($row in 36..57) {
    if (-not $oldValues.ContainsKey($row)) { continue }
    $pair = $oldValues[$row]

    # row 55-57 previously carried a shared-string note in G; clear it first
    $ws.Range("G$row").Clear()

    $ws.Range("G$row").Value = $pair[0]
    $ws.Range("H$row").Value = $pair[1]
    Set-Style $ws.Range("G$row`:H$row") $style16
}

# Row 36 & 37 keep a *new* current-run value in B (old B/C move to G/H above);
# rows 38-57 lose their old B/C entirely (cleared -> formulas fall back to 0).
$ws.Range("B36").Value = 12987
$ws.Range("C36").Clear()

$ws.Range("A37").Value = "1st Move"
$ws.Range("B37").Value = 13217
$ws.Range("C37").Clear()

foreach ($row in 38..57) {
    $ws.Range("B$row").Clear()
    $ws.Range("C$row").Clear()
}

# ---- sheetView: pane/selection moved along with the active work area ----
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("B38").Select()
